# This script applies the 26.04.2022 Slovak COVID daily-stats update:
#   1) Corrects the AgTests (F) / AgPosit (G) figures for rows 712-768
#      (these were revised upward in the source data refresh), and
#   2) Appends the AgTests/AgPosit figures that were missing for row 769,
#      plus 13 brand-new daily rows (770-782, 2022-04-13 .. 2022-04-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map column letters used below to their 1-based column index.
$colIndex = @{ "A" = 1; "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "G" = 7 }

# --- Part 1: revised AgTests / AgPosit values for existing rows 712-768 ---
# Each tuple is (Row, ColumnLetter, NewValue).
$cellUpdates = @(
    @(712, "F", 51390),
    @(715, "F", 31835),
    @(717, "F", 12568),
    @(718, "F", 17110),
    @(718, "G", 2866),
    @(721, "F", 27926),
    @(721, "G", 3149),
    @(722, "F", 27949),
    @(723, "F", 22563),
    @(726, "F", 35816),
    @(727, "F", 25203),
    @(727, "G", 2824),
    @(728, "F", 24756),
    @(728, "G", 2629),
    @(729, "F", 23296),
    @(729, "G", 2532),
    @(730, "F", 19578),
    @(730, "G", 2344),
    @(731, "F", 8673),
    @(731, "G", 1334),
    @(732, "F", 11898),
    @(732, "G", 1912),
    @(733, "F", 31765),
    @(733, "G", 3738),
    @(734, "F", 23170),
    @(734, "G", 2554),
    @(735, "F", 19346),
    @(735, "G", 2277),
    @(736, "F", 19625),
    @(736, "G", 2202),
    @(737, "F", 18526),
    @(737, "G", 2308),
    @(738, "F", 6879),
    @(738, "G", 1009),
    @(739, "F", 8706),
    @(739, "G", 1415),
    @(740, "F", 24708),
    @(740, "G", 2762),
    @(741, "F", 18983),
    @(741, "G", 1942),
    @(742, "F", 17329),
    @(742, "G", 1703),
    @(743, "F", 18072),
    @(743, "G", 1627),
    @(744, "F", 14706),
    @(744, "G", 1602),
    @(745, "F", 6210),
    @(745, "G", 936),
    @(746, "F", 8002),
    @(746, "G", 1252),
    @(747, "F", 22340),
    @(747, "G", 2384),
    @(748, "F", 16924),
    @(748, "G", 1545),
    @(749, "F", 14808),
    @(749, "G", 1482),
    @(750, "F", 15049),
    @(750, "G", 1358),
    @(751, "F", 12518),
    @(751, "G", 1384),
    @(752, "F", 4782),
    @(752, "G", 618),
    @(753, "F", 6829),
    @(753, "G", 960),
    @(754, "F", 20982),
    @(754, "G", 1945),
    @(755, "F", 13702),
    @(755, "G", 1297),
    @(756, "F", 13701),
    @(756, "G", 1074),
    @(757, "F", 13526),
    @(757, "G", 1008),
    @(758, "F", 11197),
    @(758, "G", 928),
    @(759, "F", 3873),
    @(759, "G", 387),
    @(760, "F", 5131),
    @(760, "G", 559),
    @(761, "F", 16923),
    @(761, "G", 1273),
    @(762, "F", 11265),
    @(762, "G", 791),
    @(763, "F", 10322),
    @(763, "G", 726),
    @(764, "F", 11171),
    @(764, "G", 692),
    @(765, "F", 9197),
    @(765, "G", 603),
    @(766, "F", 3218),
    @(766, "G", 252),
    @(767, "F", 4127),
    @(767, "G", 313),
    @(768, "F", 15027),
    @(768, "G", 772)
)

foreach ($u in $cellUpdates) {
    $row = $u[0]
    $col = $colIndex[$u[1]]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}

# --- Part 2: newly reported cells - the missing F769/G769 values, and
#     the full contents (columns A-G) of new rows 770-782 ---
$newRowData = @(
    @(769, "F", 10074),
    @(769, "G", 514),
    @(770, "A", 44664),
    @(770, "B", 1759793),
    @(770, "C", 8116),
    @(770, "D", 2676),
    @(770, "E", 19670),
    @(770, "F", 9105),
    @(770, "G", 425),
    @(771, "A", 44665),
    @(771, "B", 1762270),
    @(771, "C", 7148),
    @(771, "D", 2477),
    @(771, "E", 19697),
    @(771, "F", 9202),
    @(771, "G", 405),
    @(772, "A", 44666),
    @(772, "B", 1763197),
    @(772, "C", 2548),
    @(772, "D", 927),
    @(772, "E", 19721),
    @(772, "F", 2611),
    @(772, "G", 91),
    @(773, "A", 44667),
    @(773, "B", 1763995),
    @(773, "C", 3262),
    @(773, "D", 798),
    @(773, "E", 19741),
    @(773, "F", 2986),
    @(773, "G", 260),
    @(774, "A", 44668),
    @(774, "B", 1765258),
    @(774, "C", 3604),
    @(774, "D", 1263),
    @(774, "E", 19759),
    @(774, "F", 2109),
    @(774, "G", 111),
    @(775, "A", 44669),
    @(775, "B", 1765962),
    @(775, "C", 3331),
    @(775, "D", 704),
    @(775, "E", 19773),
    @(775, "F", 3227),
    @(775, "G", 201),
    @(776, "A", 44670),
    @(776, "B", 1768045),
    @(776, "C", 8971),
    @(776, "D", 2083),
    @(776, "E", 19790),
    @(776, "F", 14724),
    @(776, "G", 665),
    @(777, "A", 44671),
    @(777, "B", 1770460),
    @(777, "C", 8416),
    @(777, "D", 2415),
    @(777, "E", 19803),
    @(777, "F", 10338),
    @(777, "G", 445),
    @(778, "A", 44672),
    @(778, "B", 1772201),
    @(778, "C", 5861),
    @(778, "D", 1741),
    @(778, "E", 19817),
    @(778, "F", 8929),
    @(778, "G", 357),
    @(779, "A", 44673),
    @(779, "B", 1773653),
    @(779, "C", 6108),
    @(779, "D", 1452),
    @(779, "E", 19829),
    @(779, "F", 5804),
    @(779, "G", 255),
    @(780, "A", 44674),
    @(780, "B", 1774808),
    @(780, "C", 4910),
    @(780, "D", 1155),
    @(780, "E", 19839),
    @(780, "F", 1917),
    @(780, "G", 91),
    @(781, "A", 44675),
    @(781, "B", 1775178),
    @(781, "C", 1888),
    @(781, "D", 370),
    @(781, "E", 19852),
    @(781, "F", 1265),
    @(781, "G", 84),
    @(782, "A", 44676),
    @(782, "B", 1776576),
    @(782, "C", 7671),
    @(782, "D", 1398),
    @(782, "E", 19862),
    @(782, "F", 19),
    @(782, "G", 0)
)

foreach ($u in $newRowData) {
    $row = $u[0]
    $col = $colIndex[$u[1]]
    $val = $u[2]
    $ws.Cells.Item($row, $col).Value = $val
}
